# Updates the cryptos list with latest price / 1h volume-change figures
# (and swaps the ONDO / FirstDigitalUSD rows), matching the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D edited cells to Text format so numeric-looking strings
# (e.g. '3.509.88', '0.0785') are preserved verbatim as text, matching
# the original inlineStr cell type instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "61.414.55"
$ws.Range("D3").Value = "3.379.60"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "576.78"
$ws.Range("D6").Value = "134.96"
$ws.Range("D8").Value = "3.377.10"
$ws.Range("D13").Value = "3.952.25"
$ws.Range("D16").Value = "3.377.16"
$ws.Range("D17").Value = "25.22"
$ws.Range("D18").Value = "61.367.98"
$ws.Range("D19").Value = "14.03"
$ws.Range("D20").Value = "5.82"
$ws.Range("D22").Value = "379.83"
$ws.Range("D24").Value = "3.509.88"
$ws.Range("D26").Value = "70.71"
$ws.Range("D28").Value = "1.70"
$ws.Range("D29").Value = "7.80"
$ws.Range("D30").Value = "0.998"
$ws.Range("D32").Value = "2.18"
$ws.Range("D35").Value = "3.408.76"
$ws.Range("D36").Value = "23.42"
$ws.Range("D37").Value = "5.59"
$ws.Range("D38").Value = "6.99"
$ws.Range("D40").Value = "162.93"
$ws.Range("D41").Value = "0.0785"
$ws.Range("D42").Value = "1.24"
$ws.Range("D43").Value = "1.00"
$ws.Range("D44").Value = "4.45"
$ws.Range("D45").Value = "41.70"
$ws.Range("D47").Value = "1.63"
$ws.Range("D48").Value = "23.49"
$ws.Range("D50").Value = "23.38"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  +6.54%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("E19").Value = "  +6.31%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  +10.00%  "
$ws.Range("E28").Value = "  +19.87%  "
$ws.Range("E29").Value = "  +11.13%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("E37").Value = "  +4.95%  "
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("B42").Value = "ONDO"
$ws.Range("C42").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("E42").Value = "  +13.06%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -3.00%  "
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("E50").Value = "  +14.82%  "
$ws.Range("E51").Value = "  +4.09%  "
